$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5561593
$ws.Range("I76").Value = 18522328
$ws.Range("J76").Value = 6992.143
$ws.Range("K76").Value = 18522328
$ws.Range("L76").Value = 6992.143
$ws.Range("M76").Value = -18522013
$ws.Range("N76").Value = -7622.143
$ws.Range("H79").Value = 5561593
$ws.Range("I79").Value = 18522328
$ws.Range("J79").Value = 6992.143
$ws.Range("K79").Value = 18522328
$ws.Range("L79").Value = 6992.143
$ws.Range("M79").Value = -18521236
$ws.Range("N79").Value = -9176.143
$ws.Range("H80").Value = 4099.9287
$ws.Range("I80").Value = 1300.5714
$ws.Range("J80").Value = 6899.2856
$ws.Range("K80").Value = 3901.7142
$ws.Range("L80").Value = 20697.8568
$ws.Range("M80").Value = -2903.7142
$ws.Range("N80").Value = -22693.8568
$ws.Range("H83").Value = 4099.9287
$ws.Range("I83").Value = 1300.5714
$ws.Range("J83").Value = 6899.2856
$ws.Range("K83").Value = 11705.1426
$ws.Range("L83").Value = 62093.5704
$ws.Range("M83").Value = -6713.142600000001
$ws.Range("N83").Value = -72077.5704
$ws.Range("H86").Value = 8494.440000000001
$ws.Range("I86").Value = 6167.3335
$ws.Range("J86").Value = 10642.538
$ws.Range("K86").Value = 6167.3335
$ws.Range("L86").Value = 10642.538
$ws.Range("M86").Value = -5044.3335
$ws.Range("N86").Value = -12888.538
$ws.Range("H89").Value = 8494.440000000001
$ws.Range("I89").Value = 6167.3335
$ws.Range("J89").Value = 10642.538
$ws.Range("K89").Value = 30836.6675
$ws.Range("L89").Value = 53212.69
$ws.Range("M89").Value = -25220.6675
$ws.Range("N89").Value = -64444.69
$ws.Range("H92").Value = 858.6579
$ws.Range("I92").Value = 274.10345
$ws.Range("J92").Value = 2742.2222
$ws.Range("K92").Value = 274.10345
$ws.Range("L92").Value = 2742.2222
$ws.Range("M92").Value = 973.8965499999999
$ws.Range("N92").Value = -5238.2222
$ws.Range("H94").Value = 195.6
$ws.Range("I94").Value = 195.6
$ws.Range("K94").Value = 195.6
$ws.Range("M94").Value = 255.4
$ws.Range("H96").Value = 254.71428
$ws.Range("I96").Value = 398.2857
$ws.Range("K96").Value = 1194.8571
$ws.Range("M96").Value = 178.1428999999998
$ws.Range("H100").Value = 3874.75
$ws.Range("I100").Value = 4000
$ws.Range("J100").Value = 3749.5
$ws.Range("K100").Value = 4000
$ws.Range("L100").Value = 3749.5
$ws.Range("M100").Value = -3459
$ws.Range("N100").Value = -4831.5
$ws.Range("H116").Value = 8481.6
$ws.Range("I116").Value = 5749.5
$ws.Range("J116").Value = 9164.625
$ws.Range("K116").Value = 5749.5
$ws.Range("L116").Value = 9164.625
$ws.Range("M116").Value = -2307.5
$ws.Range("N116").Value = -16048.625
$ws.Range("H137").Value = 53095.6
$ws.Range("I137").Value = 67716.59
$ws.Range("K137").Value = 203149.77
$ws.Range("M137").Value = -200599.77

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5658829
$ws.Range("I2").Value = 8082792
$ws.Range("K2").Value = 8082792
$ws.Range("M2").Value = -8082679
$ws.Range("H32").Value = 4200.27
$ws.Range("I32").Value = 2175.378
$ws.Range("J32").Value = 13424.777
$ws.Range("K32").Value = 2175.378
$ws.Range("L32").Value = 13424.777
$ws.Range("M32").Value = -1888.378
$ws.Range("N32").Value = -13998.777
$ws.Range("H45").Value = 9595447
$ws.Range("J45").Value = 7747.75
$ws.Range("L45").Value = 7747.75
$ws.Range("N45").Value = -8501.75
$ws.Range("H61").Value = 2996.4167
$ws.Range("I61").Value = 2384.6667
$ws.Range("K61").Value = 2384.6667
$ws.Range("M61").Value = -2172.6667
$ws.Range("H64").Value = 24996.25
$ws.Range("J64").Value = 24996.25
$ws.Range("L64").Value = 24996.25
$ws.Range("N64").Value = -25492.25
$ws.Range("H67").Value = 24996.25
$ws.Range("J67").Value = 24996.25
$ws.Range("L67").Value = 24996.25
$ws.Range("N67").Value = -26712.25
$ws.Range("H97").Value = 984104.25
$ws.Range("I97").Value = 1475291
$ws.Range("J97").Value = 1730.8182
$ws.Range("K97").Value = 1475291
$ws.Range("L97").Value = 1730.8182
$ws.Range("M97").Value = -1474795
$ws.Range("N97").Value = -2722.8182
$ws.Range("H116").Value = 5658829
$ws.Range("I116").Value = 8082792
$ws.Range("K116").Value = 8082792
$ws.Range("M116").Value = -8080498
$ws.Range("H136").Value = 2996.4167
$ws.Range("I136").Value = 2384.6667
$ws.Range("K136").Value = 7154.000100000001
$ws.Range("M136").Value = -4604.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5658829
$ws.Range("I3").Value = 8082792
$ws.Range("K3").Value = 8082792
$ws.Range("M3").Value = -8082678
$ws.Range("H94").Value = 2329139.8
$ws.Range("I94").Value = 3126522.5
$ws.Range("J94").Value = 9481.091
$ws.Range("K94").Value = 3126522.5
$ws.Range("L94").Value = 9481.091
$ws.Range("M94").Value = -3126071.5
$ws.Range("N94").Value = -10383.091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 8026.8335
$ws.Range("J14").Value = 8026.8335
$ws.Range("L14").Value = 8026.8335
$ws.Range("N14").Value = -8366.833500000001
$ws.Range("H16").Value = 1900.4286
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H38").Value = 27500
$ws.Range("I38").Value = 45000
$ws.Range("K38").Value = 45000
$ws.Range("M38").Value = -44623
$ws.Range("H46").Value = 27500
$ws.Range("I46").Value = 45000
$ws.Range("K46").Value = 45000
$ws.Range("M46").Value = -44789
$ws.Range("H58").Value = 9709.6
$ws.Range("J58").Value = 4945
$ws.Range("L58").Value = 4945
$ws.Range("N58").Value = -5351
$ws.Range("H105").Value = 1999.5
$ws.Range("I105").Value = 1999.5
$ws.Range("K105").Value = 1999.5
$ws.Range("M105").Value = -252.5
$ws.Range("H107").Value = 83337170
$ws.Range("J107").Value = 200003940
$ws.Range("L107").Value = 200003940
$ws.Range("N107").Value = -200007780
$ws.Range("H113").Value = 1900.4286
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 102388.48
$ws.Range("I132").Value = 114055.336
$ws.Range("K132").Value = 342166.008
$ws.Range("M132").Value = -339636.008
$ws.Range("H136").Value = 9709.6
$ws.Range("J136").Value = 4945
$ws.Range("L136").Value = 14835
$ws.Range("N136").Value = -19935

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 50807.1
$ws.Range("I5").Value = 747.6429000000001
$ws.Range("J5").Value = 167612.5
$ws.Range("K5").Value = 2242.9287
$ws.Range("L5").Value = 502837.5
$ws.Range("M5").Value = -2130.9287
$ws.Range("N5").Value = -503061.5
$ws.Range("H131").Value = 15436156
$ws.Range("I131").Value = 11908008
$ws.Range("J131").Value = 16671008
$ws.Range("K131").Value = 35724024
$ws.Range("L131").Value = 50013024
$ws.Range("M131").Value = -35718984
$ws.Range("N131").Value = -50023104
$ws.Range("H135").Value = 50807.1
$ws.Range("I135").Value = 747.6429000000001
$ws.Range("J135").Value = 167612.5
$ws.Range("K135").Value = 6728.7861
$ws.Range("L135").Value = 1508512.5
$ws.Range("M135").Value = -4193.7861
$ws.Range("N135").Value = -1513582.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H97").Value = 2977579.2
$ws.Range("I97").Value = 4762947
$ws.Range("K97").Value = 4762947
$ws.Range("M97").Value = -4762451
$ws.Range("H126").Value = 3791747
$ws.Range("I126").Value = 4548153.5
$ws.Range("J126").Value = 3476577.5
$ws.Range("K126").Value = 13644460.5
$ws.Range("L126").Value = 10429732.5
$ws.Range("M126").Value = -13641990.5
$ws.Range("N126").Value = -10434672.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 10424748
$ws.Range("I93").Value = 15153484
$ws.Range("J93").Value = 21528.9
$ws.Range("K93").Value = 15153484
$ws.Range("L93").Value = 21528.9
$ws.Range("M93").Value = -15152236
$ws.Range("N93").Value = -24024.9
$ws.Range("H100").Value = 3413.087
$ws.Range("I100").Value = 2889
$ws.Range("K100").Value = 2889
$ws.Range("M100").Value = -2348

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4737
$ws.Range("I122").Value = 4149.778
$ws.Range("K122").Value = 12449.334
$ws.Range("M122").Value = -9999.334000000001
$ws.Range("H126").Value = 2856.1428
$ws.Range("I126").Value = 3174.6667
$ws.Range("K126").Value = 9524.000100000001
$ws.Range("M126").Value = -7054.000100000001
